$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the Actual Result / Pass-Fail columns for the two test cases
$ws.Range("H3").Value = "Error sign. Email already registered with another account"
$ws.Range("H2").Value = "Successfully Navigated to Home page"
$ws.Range("I2").Value = "Pass"
$ws.Range("I3").Value = "Pass"

# Clear the stray empty cells left in rows 4-5 (C4:D5)
$ws.Range("C4:D5").ClearContents()

# Update the active selection to I3, matching the author's final cursor position
$ws.Range("I3").Select()
